$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Content fix: row 17 (category favourited in UI test) pointed at the
#     wrong tab name -- correct "seller tab" -> "category tab".
$ws.Range("C17").Value = "The category is saved in the category tab"

# --- Highlight the newly-implemented "UI favourites tests" rows (14-23)
#     using the same green fill already used further down the sheet for
#     completed rows (style ids 8/9 reuse the existing themed fill).
$ws.Range("A28").Copy()
$ws.Range("A14:A23").PasteSpecial(-4122)
$ws.Range("C28").Copy()
$ws.Range("C14:C23").PasteSpecial(-4122)

$ws.Range("B13").Copy()
$ws.Range("B14:B20").PasteSpecial(-4122)
$ws.Range("B28").Copy()
$ws.Range("B21:B23").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# --- Move the view back up to the top of the newly finished section and
#     select the last completed row.
$ws.Activate()
$excel.ActiveWindow.Zoom = 100
$ws.Range("A23:C23").Select()
